$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "about getting access (try to do this by the end of February).",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "about getting access.",
    2
)
